$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Helper: find the 1-based row index whose first cell ("Nom de balise")
# equals the given tag name (cell text carries a trailing cell-mark, so
# strip CR/BEL/vertical-tab before comparing).
function Get-RowIndexByTag($table, $tag) {
    for ($i = 1; $i -le $table.Rows.Count; $i++) {
        $txt = $table.Cell($i, 1).Range.Text
        $txt = $txt -replace "[\x07\r\x0B]", ""
        if ($txt -eq $tag) {
            return $i
        }
    }
    return -1
}

# --- Change 1: "objectField" row, Format column: object -> technicalObject ---
$objectFieldRow = Get-RowIndexByTag $t "objectField"
$t.Cell($objectFieldRow, 3).Range.Text = "technicalObject"

# --- Change 2: insert three new rows right after the "arrayField" row ---
$arrayFieldRow = Get-RowIndexByTag $t "arrayField"
$insertPos = $arrayFieldRow + 1

$insertBefore = $t.Rows.Item($insertPos)
$row1 = $t.Rows.Add($insertBefore)
$row1.Cells.Item(1).Range.Text = "enumArrayField"
$row1.Cells.Item(2).Range.Text = "Array of enumerations"
$row1.Cells.Item(3).Range.Text = "string" + [char]11 + "(ENUM: REPORT, UPDATE, CANCEL, ACK, ERROR)"
$row1.Cells.Item(4).Range.Text = "0..n"
$row1.Cells.Item(5).Range.Text = "This is an array of enumerations"

$insertBefore = $t.Rows.Item($insertPos + 1)
$row2 = $t.Rows.Add($insertBefore)
$row2.Cells.Item(1).Range.Text = "requiredArray"
$row2.Cells.Item(2).Range.Text = "Required array"
$row2.Cells.Item(3).Range.Text = "string"
$row2.Cells.Item(4).Range.Text = "1..n"
$row2.Cells.Item(5).Range.Text = "This array is required"

$insertBefore = $t.Rows.Item($insertPos + 2)
$row3 = $t.Rows.Add($insertBefore)
$row3.Cells.Item(1).Range.Text = "arrayWithMaxLength"
$row3.Cells.Item(2).Range.Text = "Array with maximum length"
$row3.Cells.Item(3).Range.Text = "string"
$row3.Cells.Item(4).Range.Text = "0..5"
$row3.Cells.Item(5).Range.Text = "This is an array with a maximum length"

# --- Change 3: "objectLevel1" row, Format column: object -> levelOneData ---
$objectLevel1Row = Get-RowIndexByTag $t "objectLevel1"
$t.Cell($objectLevel1Row, 3).Range.Text = "levelOneData"
